# Auto-generated Excel COM-interop script
# Updates numeric cell values across multiple worksheets
# as part of a scheduled data-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1016.5455
$ws.Range("I129").Value = 821.3333
$ws.Range("K129").Value = 2463.9999
$ws.Range("M129").Value = 2536.0001
$ws.Range("H131").Value = 1126.4
$ws.Range("I131").Value = 1126.4
$ws.Range("K131").Value = 3379.2
$ws.Range("M131").Value = 1660.8
$ws.Range("H138").Value = 3770.3333
$ws.Range("I138").Value = 2375.077
$ws.Range("J138").Value = 7398
$ws.Range("K138").Value = 7125.231000000001
$ws.Range("L138").Value = 22194
$ws.Range("M138").Value = -1985.231000000001
$ws.Range("N138").Value = -32474

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 116624.1
$ws.Range("I32").Value = 121721.86
$ws.Range("K32").Value = 121721.86
$ws.Range("M32").Value = -121434.86
$ws.Range("H74").Value = 1075.2609
$ws.Range("I74").Value = 879.9729599999999
$ws.Range("J74").Value = 1878.1111
$ws.Range("K74").Value = 879.9729599999999
$ws.Range("L74").Value = 1878.1111
$ws.Range("M74").Value = -5.972959999999944
$ws.Range("N74").Value = -3626.1111
$ws.Range("H77").Value = 1075.2609
$ws.Range("I77").Value = 879.9729599999999
$ws.Range("J77").Value = 1878.1111
$ws.Range("K77").Value = 4399.864799999999
$ws.Range("L77").Value = 9390.5555
$ws.Range("M77").Value = -31.86479999999938
$ws.Range("N77").Value = -18126.5555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3637.532
$ws.Range("I31").Value = 3140.6316
$ws.Range("J31").Value = 3974.7144
$ws.Range("K31").Value = 3140.6316
$ws.Range("L31").Value = 3974.7144
$ws.Range("M31").Value = -2845.6316
$ws.Range("N31").Value = -4564.7144
$ws.Range("H34").Value = 3637.532
$ws.Range("I34").Value = 3140.6316
$ws.Range("J34").Value = 3974.7144
$ws.Range("K34").Value = 3140.6316
$ws.Range("L34").Value = 3974.7144
$ws.Range("M34").Value = -2938.6316
$ws.Range("N34").Value = -4378.7144
$ws.Range("H58").Value = 1555.3636
$ws.Range("I58").Value = 1589.8889
$ws.Range("J58").Value = 1400
$ws.Range("K58").Value = 1589.8889
$ws.Range("L58").Value = 1400
$ws.Range("M58").Value = -1386.8889
$ws.Range("N58").Value = -1806
$ws.Range("H132").Value = 2755.2222
$ws.Range("I132").Value = 2643.3635
$ws.Range("J132").Value = 3985.6667
$ws.Range("K132").Value = 7930.0905
$ws.Range("L132").Value = 11957.0001
$ws.Range("M132").Value = -5400.0905
$ws.Range("N132").Value = -17017.0001
$ws.Range("H136").Value = 1555.3636
$ws.Range("I136").Value = 1589.8889
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 4769.6667
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -2219.6667
$ws.Range("N136").Value = -9300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 126074.5
$ws.Range("I34").Value = 500249.5
$ws.Range("J34").Value = 1349.5
$ws.Range("K34").Value = 1500748.5
$ws.Range("L34").Value = 4048.5
$ws.Range("M34").Value = -1500664.5
$ws.Range("N34").Value = -4216.5
$ws.Range("H46").Value = 123957.11
$ws.Range("I46").Value = 999999
$ws.Range("J46").Value = 14451.875
$ws.Range("K46").Value = 2999997
$ws.Range("L46").Value = 43355.625
$ws.Range("M46").Value = -2999906
$ws.Range("N46").Value = -43537.625
$ws.Range("H50").Value = 335332.66
$ws.Range("I50").Value = 999999
$ws.Range("J50").Value = 2999.5
$ws.Range("K50").Value = 2999997
$ws.Range("L50").Value = 8998.5
$ws.Range("M50").Value = -2999516
$ws.Range("N50").Value = -9960.5
$ws.Range("H53").Value = 335332.66
$ws.Range("I53").Value = 999999
$ws.Range("J53").Value = 2999.5
$ws.Range("K53").Value = 2999997
$ws.Range("L53").Value = 8998.5
$ws.Range("M53").Value = -2999516
$ws.Range("N53").Value = -9960.5
$ws.Range("H55").Value = 5482369
$ws.Range("I55").Value = 168931.33
$ws.Range("J55").Value = 7357700
$ws.Range("K55").Value = 506793.99
$ws.Range("L55").Value = 22073100
$ws.Range("M55").Value = -506616.99
$ws.Range("N55").Value = -22073454
$ws.Range("H68").Value = 3575610.8
$ws.Range("I68").Value = 3999.75
$ws.Range("J68").Value = 8337759
$ws.Range("K68").Value = 11999.25
$ws.Range("L68").Value = 25013277
$ws.Range("M68").Value = -11188.25
$ws.Range("N68").Value = -25014899
$ws.Range("H71").Value = 3575610.8
$ws.Range("I71").Value = 3999.75
$ws.Range("J71").Value = 8337759
$ws.Range("K71").Value = 35997.75
$ws.Range("L71").Value = 75039831
$ws.Range("M71").Value = -31941.75
$ws.Range("N71").Value = -75047943
$ws.Range("H94").Value = 122077.445
$ws.Range("I94").Value = 501199
$ws.Range("J94").Value = 13757
$ws.Range("K94").Value = 1503597
$ws.Range("L94").Value = 41271
$ws.Range("M94").Value = -1502921
$ws.Range("N94").Value = -42623
$ws.Range("H104").Value = 103729.8
$ws.Range("I104").Value = 251474.75
$ws.Range("J104").Value = 5233.1665
$ws.Range("K104").Value = 754424.25
$ws.Range("L104").Value = 15699.4995
$ws.Range("M104").Value = -751803.25
$ws.Range("N104").Value = -20941.4995
$ws.Range("H107").Value = 2285.5833
$ws.Range("I107").Value = 1483.75
$ws.Range("K107").Value = 4451.25
$ws.Range("M107").Value = -2531.25
$ws.Range("H113").Value = 465.70834
$ws.Range("J113").Value = 612.5
$ws.Range("L113").Value = 1837.5
$ws.Range("N113").Value = -6177.5
$ws.Range("H118").Value = 204397.8
$ws.Range("J118").Value = 5497.5
$ws.Range("L118").Value = 16492.5
$ws.Range("N118").Value = -18978.5
$ws.Range("H121").Value = 15243677
$ws.Range("I121").Value = 25718594
$ws.Range("J121").Value = 113239.22
$ws.Range("K121").Value = 77155782
$ws.Range("L121").Value = 339717.66
$ws.Range("M121").Value = -77154472
$ws.Range("N121").Value = -342337.66

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58832810
$ws.Range("I80").Value = 125014460
$ws.Range("J80").Value = 4665.5557
$ws.Range("K80").Value = 125014460
$ws.Range("L80").Value = 4665.5557
$ws.Range("M80").Value = -125013462
$ws.Range("N80").Value = -6661.5557
$ws.Range("H83").Value = 58832810
$ws.Range("I83").Value = 125014460
$ws.Range("J83").Value = 4665.5557
$ws.Range("K83").Value = 625072300
$ws.Range("L83").Value = 23327.7785
$ws.Range("M83").Value = -625067308
$ws.Range("N83").Value = -33311.7785
$ws.Range("H113").Value = 8327.895
$ws.Range("I113").Value = 8764.25
$ws.Range("J113").Value = 6000.6665
$ws.Range("K113").Value = 8764.25
$ws.Range("L113").Value = 6000.6665
$ws.Range("M113").Value = -6594.25
$ws.Range("N113").Value = -10340.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 88500
$ws.Range("I7").Value = 130000
$ws.Range("J7").Value = 5500
$ws.Range("K7").Value = 130000
$ws.Range("L7").Value = 5500
$ws.Range("M7").Value = -129888
$ws.Range("N7").Value = -5724
$ws.Range("H63").Value = 24938.75
$ws.Range("J63").Value = 24938.75
$ws.Range("L63").Value = 24938.75
$ws.Range("N63").Value = -26436.75
$ws.Range("H66").Value = 24938.75
$ws.Range("J66").Value = 24938.75
$ws.Range("L66").Value = 74816.25
$ws.Range("N66").Value = -82304.25
$ws.Range("H93").Value = 2025.1111
$ws.Range("I93").Value = 2083.3333
$ws.Range("J93").Value = 1996
$ws.Range("K93").Value = 2083.3333
$ws.Range("L93").Value = 1996
$ws.Range("M93").Value = -835.3332999999998
$ws.Range("N93").Value = -4492
$ws.Range("H126").Value = 88500
$ws.Range("I126").Value = 130000
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 390000
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -387530
$ws.Range("N126").Value = -21440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 12914.429
$ws.Range("I4").Value = 14667.333
$ws.Range("J4").Value = 11599.75
$ws.Range("K4").Value = 14667.333
$ws.Range("L4").Value = 11599.75
$ws.Range("N4").Value = -11825.75
$ws.Range("M4").Value = -14554.333
$ws.Range("H126").Value = 1493.909
$ws.Range("I126").Value = 1384.7778
$ws.Range("J126").Value = 1985
$ws.Range("K126").Value = 4154.3334
$ws.Range("L126").Value = 5955
$ws.Range("M126").Value = -1684.3334
$ws.Range("N126").Value = -10895
